$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 150, shifting existing rows 150:173 down to 151:174.
$ws.Rows("150:150").Insert()

# Populate the newly inserted row 150 with the new weekly record.
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 44491
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = 100112008
$ws.Cells.Item(150, 7).Value = "Coliflor"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 200
$ws.Cells.Item(150, 11).Value = 700
$ws.Cells.Item(150, 12).Value = 800
$ws.Cells.Item(150, 13).Value = 750
$ws.Cells.Item(150, 14).Value = "$/unidad"
$ws.Cells.Item(150, 15).Value = "Región Metropolitana"
$ws.Cells.Item(150, 16).Value = 750
$ws.Cells.Item(150, 17).Value = 1
$ws.Cells.Item(150, 18).Value = "Hortaliza"
